$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 4 down to row 5 so new cells inherit number formats / styles
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row 5 values
$ws.Range("A5").Value = 43986
$ws.Range("B5").Value = 105680
$ws.Range("C5").Value = 161724
$ws.Range("D5").Value = 46659
$ws.Range("E5").Value = 12545
$ws.Range("F5").Value = 34.235427706283119
$ws.Range("G5").Value = 36180
$ws.Range("H5").Value = 3405
$ws.Range("I5").Value = 3587
$ws.Range("J5").Value = 314063

# Update the selected cell to match the new view state
$ws.Range("C8").Select()
